$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns: P1=14, Q1=15
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Give the new header cells the same formatting as the existing header cell O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats

# Update existing columns I, K, M, O for rows 2-25 (values swap 1<->2),
# and populate the two newly added columns P and Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # column I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # column K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # column M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # column O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # column P (new) -> 2
    $ws.Cells.Item($r, 17).Value = 2   # column Q (new) -> 2
}
